$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 911, shifting existing rows 911..1011 down to 912..1012.
$ws.Rows(911).Insert()

# Populate the newly inserted row 911 with the new weekly data point.
$ws.Cells.Item(911, 1).Value = 4
$ws.Cells.Item(911, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(911, 3).Value = "Los Lagos"
$ws.Cells.Item(911, 4).Value = 45194
$ws.Cells.Item(911, 5).Value = 10
$ws.Cells.Item(911, 6).Value = 100112004
$ws.Cells.Item(911, 7).Value = "Cebolla"
$ws.Cells.Item(911, 8).Value = "Sin especificar"
$ws.Cells.Item(911, 9).Value = "1a (guarda)"
$ws.Cells.Item(911, 10).Value = 500
$ws.Cells.Item(911, 11).Value = 18000
$ws.Cells.Item(911, 12).Value = 18000
$ws.Cells.Item(911, 13).Value = 18000
$ws.Cells.Item(911, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(911, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(911, 16).Value = 1000
$ws.Cells.Item(911, 17).Value = 18
$ws.Cells.Item(911, 18).Value = "Hortaliza"
